$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the header labels:
#    "<Name>_old" -> "<Name>_FV2310"   (columns A:J)
#    "<Name>_new" -> "<Name>_FV2404"   (columns L:U)
#    column K ("diff") is left untouched
# ---------------------------------------------------------------------
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# ---------------------------------------------------------------------
# 2) Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3) Turn the data range into an Excel table (ListObject) that exposes
#    the freshly renamed headers as its column names.  The header row
#    already carries manual formatting (bold font, grey fill, border);
#    stash it away, reset the range to the default style before the
#    table is created (otherwise the engine would bake the manual
#    formatting into a dedicated header dxf), and restore the
#    look afterwards.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$headerRange.Style = "Normal"

$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.Clear()
$tbl.TableStyle = ""
